# Update the "education" worksheet with the new set of degrees/records,
# apply wrap-text/left-top alignment formatting to the data rows, set the
# resulting (2-line) row heights, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New table content (what / when / with / where) -----------------------
# Row 2
$ws.Range("A2").Value2 = "PhD - Neurosciences"
$ws.Range("B2").Value2 = 2018
$ws.Range("C2").Value2 = "\href{https://www.uv.es/uvweb/universidad/es/universidad-valencia-1285845048380.html}{Universidad de Valencia}"
$ws.Range("D2").Value2 = "Valencia, España"

# Row 3
$ws.Range("A3").Value2 = "Master's Degree in Basic and Applied Neurosciences "
$ws.Range("B3").Value2 = 2012
$ws.Range("C3").Value2 = "\href{https://www.uv.es/uvweb/universidad/es/universidad-valencia-1285845048380.html}{Universidad de Valencia}"
$ws.Range("D3").Value2 = "Valencia, España"

# Row 4
$ws.Range("A4").Value2 = "Psychology "
$ws.Range("B4").Value2 = 2007
$ws.Range("C4").Value2 = "\href{https://www.ucatolica.edu.co/portal/Pregrado/psicologia/}{Universidad Cátolica de Colombia}"
$ws.Range("D4").Value2 = "Bogotá, Colombia"

# --- Formatting: left/top aligned, wrapped text for the data rows ---------
# Build the combined alignment format once on a scratch cell (so only a
# single new cell style is produced), then copy just the formatting onto
# the data range and discard the scratch cell.
$tmpl = $ws.Cells.Item(20, 20)
$tmpl.HorizontalAlignment = -4131   # xlLeft
$tmpl.VerticalAlignment = -4160     # xlTop
$tmpl.WrapText = $true

$tmpl.Copy()
$ws.Range("A2:D4").PasteSpecial(-4122)   # xlPasteFormats
$tmpl.Clear()

# Two lines of wrapped text at the sheet's (14.4pt) default row height.
$ws.Rows.Item(2).RowHeight = 28.8
$ws.Rows.Item(3).RowHeight = 28.8
$ws.Rows.Item(4).RowHeight = 28.8

# --- Selection change -------------------------------------------------------
$ws.Range("C10").Select()
